# TC03_Canine_Filter_Breed-AusShephd.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab query (row 2 / cell B2 on sheet "startup") included an extra
# `coalesce(co.cohort_description, '') AS `Cohort`` column that the fix
# removes from the returned query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
    "WHERE demo.breed  IN ['Australian Shepherd']`n" +
    "MATCH (c)<--(diag:diagnosis)`n" +
    "OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
    "WITH DISTINCT c, s, demo, diag, co`n" +
    "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
    "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
    "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
    "        coalesce(demo.breed, '') AS Breed ,`n" +
    "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
    "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
    "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
    "        coalesce(demo.sex, '') AS Sex ,`n" +
    "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
    "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
    "        coalesce(diag.best_response, '') AS ``Response to Treatment``" + "`n"

$ws.Range("B2").Value = $newCasesQuery

# Row heights shrink (content of B2 lost a line, and the other rows' wrapped
# text reflows slightly) -- match the committed values exactly.
$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# View / window state recorded with the fix.
$win = $excel.ActiveWindow
$win.Zoom = 160
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("B2").Select()

# Best-effort: the author's Excel top-level window also moved/resized when
# the fix was saved. Not all hosts expose writable window geometry, so this
# is wrapped defensively and ignored if unsupported.
try { $win.WindowState = -4143 } catch {}
try { $win.Left = 28680 } catch {}
try { $win.Top = -105 } catch {}
try { $win.Width = 29040 } catch {}
try { $win.Height = 15840 } catch {}
